$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = "kjopwej;lk34rq3[]p\]=34\=5p34=\ol54"
$ws.Range("F12").Value = "]3j[4i]lo[o-]ik23]"
$ws.Range("F13").Value = "p4-p0oj2p3lm,e;w.a]\d=[o-"
$ws.Range("F14").Value = "pio0a9ihuihzxbJ;'["
$ws.Range("F16").Value = "o-0i903ui2heknwqml;'"
$ws.Range("F17").Value = "][=p-o0i9u8yhjkm,l;'[ptdtresw45r67t98u0i-o\]["
$ws.Range("F18").Value = "'?>:,mkjhbgtfdr567t890-[-09iu8tdxfvbnjkl;'[]9807ytrdfxc"

$ws.Range("F18").Select()
